{"js": "const pairs = [\n  [\"2025-07-21 Monday\", \"2025-07-22 Tuesday\"],\n  [\"596\u00f78=74, 4\", \"326\u00f72=163, 0\"],\n  [\"339\u00f77=48, 3\", \"395\u00f73=131, 2\"],\n  [\"901\u00f75=180, 1\", \"795\u00f74=198, 3\"],\n  [\"226\u00f74=56, 2\", \"198\u00f74=49, 2\"],\n  [\"462\u00f72=231, 0\", \"708\u00f75=141, 3\"],\n  [\"634\u00f73=211, 1\", \"240\u00f76=40, 0\"],\n  [\"125\u00f72=62, 1\", \"838\u00f73=279, 1\"],\n  [\"266\u00f76=44, 2\", \"991\u00f79=110, 1\"],\n  [\"255\u00f77=36, 3\", \"295\u00f72=147, 1\"],\n  [\"384\u00f73=128, 0\", \"604\u00f73=201, 1\"],\n  [\"546\u00f73=182, 0\", \"472\u00f74=118, 0\"],\n  [\"941\u00f76=156, 5\", \"876\u00f79=97, 3\"],\n  [\"275\u00f76=45, 5\", \"676\u00f74=169, 0\"],\n  [\"963\u00f78=120, 3\", \"612\u00f73=204, 0\"],\n  [\"350\u00f74=87, 2\", \"248\u00f76=41, 2\"],\n  [\"761\u00f77=108, 5\", \"559\u00f78=69, 7\"],\n  [\"311\u00f77=44, 3\", \"551\u00f76=91, 5\"],\n  [\"676\u00f78=84, 4\", \"582\u00f78=72, 6\"],\n  [\"262\u00f79=29, 1\", \"162\u00f75=32, 2\"],\n  [\"526\u00f76=87, 4\", \"313\u00f74=78, 1\"],\n  [\"869\u00f79=96, 5\", \"148\u00f75=29, 3\"],\n  [\"233\u00f76=38, 5\", \"628\u00f73=209, 1\"],\n  [\"263\u00f75=52, 3\", \"823\u00f76=137, 1\"],\n  [\"839\u00f73=279, 2\", \"143\u00f78=17, 7\"],\n  [\"310\u00f73=103, 1\", \"627\u00f76=104, 3\"],\n];\n\nfor (const [before, after] of pairs) {\n  const results = context.document.body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(after, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2025-07-21 Monday\", \"2025-07-22 Tuesday\"),\n  @(\"596\u00f78=74, 4\", \"326\u00f72=163, 0\"),\n  @(\"339\u00f77=48, 3\", \"395\u00f73=131, 2\"),\n  @(\"901\u00f75=180, 1\", \"795\u00f74=198, 3\"),\n  @(\"226\u00f74=56, 2\", \"198\u00f74=49, 2\"),\n  @(\"462\u00f72=231, 0\", \"708\u00f75=141, 3\"),\n  @(\"634\u00f73=211, 1\", \"240\u00f76=40, 0\"),\n  @(\"125\u00f72=62, 1\", \"838\u00f73=279, 1\"),\n  @(\"266\u00f76=44, 2\", \"991\u00f79=110, 1\"),\n  @(\"255\u00f77=36, 3\", \"295\u00f72=147, 1\"),\n  @(\"384\u00f73=128, 0\", \"604\u00f73=201, 1\"),\n  @(\"546\u00f73=182, 0\", \"472\u00f74=118, 0\"),\n  @(\"941\u00f76=156, 5\", \"876\u00f79=97, 3\"),\n  @(\"275\u00f76=45, 5\", \"676\u00f74=169, 0\"),\n  @(\"963\u00f78=120, 3\", \"612\u00f73=204, 0\"),\n  @(\"350\u00f74=87, 2\", \"248\u00f76=41, 2\"),\n  @(\"761\u00f77=108, 5\", \"559\u00f78=69, 7\"),\n  @(\"311\u00f77=44, 3\", \"551\u00f76=91, 5\"),\n  @(\"676\u00f78=84, 4\", \"582\u00f78=72, 6\"),\n  @(\"262\u00f79=29, 1\", \"162\u00f75=32, 2\"),\n  @(\"526\u00f76=87, 4\", \"313\u00f74=78, 1\"),\n  @(\"869\u00f79=96, 5\", \"148\u00f75=29, 3\"),\n  @(\"233\u00f76=38, 5\", \"628\u00f73=209, 1\"),\n  @(\"263\u00f75=52, 3\", \"823\u00f76=137, 1\"),\n  @(\"839\u00f73=279, 2\", \"143\u00f78=17, 7\"),\n  @(\"310\u00f73=103, 1\", \"627\u00f76=104, 3\"),\n)\n\nforeach ($pair in $pairs) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Execute($pair[0], $false, $true, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n}"}
